# Update yearly.xlsx "Overview" sheet:
#  - shift the fiscal-year column headers forward by one year
#    (1396/12..1400/12 -> 1397/12..1401/12)
#  - refresh the underlying database values (read_price algorithm update)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column headers (row 8 and row 24) ----
$ws.Range("E8").Value  = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value  = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value  = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value  = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value  = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# ---- هزینه حمل و نقل و انتقال (row 10) ----
$ws.Range("E10").Value = 4247
$ws.Range("F10").Value = 5036
$ws.Range("G10").Value = 11594
$ws.Range("H10").Value = 143287
$ws.Range("I10").Value = 372932

# ---- هزینه انرژی (آب، برق، گاز و سوخت) (row 16) ----
$ws.Range("E16").Value = 887
$ws.Range("F16").Value = 2082
$ws.Range("G16").Value = 3336
$ws.Range("H16").Value = 4083
$ws.Range("I16").Value = 11414

# ---- هزینه استهلاک (row 17) ----
$ws.Range("E17").Value = 21101
$ws.Range("F17").Value = 26858
$ws.Range("G17").Value = 34027
$ws.Range("H17").Value = 50150
$ws.Range("I17").Value = 116851

# ---- هزینه مطالبات مشکوک الوصول (row 19) ----
$ws.Range("E19").Value = 60284
$ws.Range("F19").Value = 124915
$ws.Range("G19").Value = 215906
$ws.Range("H19").Value = 183932
$ws.Range("I19").Value = 483701

# ---- جمع - total (row 20) ----
$ws.Range("E20").Value = 86519
$ws.Range("F20").Value = 158891
$ws.Range("G20").Value = 264863
$ws.Range("H20").Value = 381452
$ws.Range("I20").Value = 984898

# ---- تعداد پرسنل غیر تولیدی شرکت (row 26) ----
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = 39
$ws.Range("I26").Value = 41

# ---- تعداد پرسنل تولیدی شرکت (row 27) ----
$ws.Range("E27").Value = 125
$ws.Range("F27").Value = 127
$ws.Range("G27").Value = 125
$ws.Range("H27").Value = 162
$ws.Range("I27").Value = 167

Write-Output "applied yearly database update"
